# Append the daily Covid-19 data rows (2020-05-25 .. 2020-07-02) to the
# "Covid-19 podatki" sheet / "Tabela1" table, following the same pattern
# the publisher's daily bot uses: copy the previous row's formatting down
# into the new row, then overwrite the 10 column values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
    @(43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0),
    @(43978, 77210, 631, 1473, 2, 7, 2, 1, 108, 0),
    @(43979, 77916, 706, 1473, 0, 7, 2, 0, 108, 0),
    @(43980, 78529, 613, 1473, 0, 7, 2, 0, 108, 0),
    @(43981, 78793, 264, 1473, 0, 6, 2, 1, 108, 0),
    @(43982, 79039, 246, 1473, 0, 5, 1, 0, 109, 1),
    @(43983, 79698, 659, 1475, 2, 5, 1, 0, 109, 0),
    @(43984, 80505, 807, 1477, 2, 5, 0, 0, 109, 0),
    @(43985, 81333, 828, 1477, 0, 5, 0, 0, 109, 0),
    @(43986, 82161, 828, 1479, 2, 6, 0, 0, 109, 0),
    @(43987, 82876, 715, 1484, 5, 6, 0, 0, 109, 0),
    @(43988, 83105, 229, 1485, 1, 5, 0, 1, 109, 0),
    @(43989, 83316, 211, 1485, 0, 5, 0, 0, 109, 0),
    @(43990, 84130, 814, 1486, 1, 6, 0, 0, 109, 0),
    @(43991, 84868, 738, 1488, 2, 6, 0, 0, 109, 0),
    @(43992, 85626, 758, 1488, 0, 6, 0, 0, 109, 0),
    @(43993, 86328, 702, 1490, 2, 6, 0, 0, 109, 0),
    @(43994, 87095, 767, 1492, 2, 6, 0, 0, 109, 0),
    @(43995, 87386, 291, 1495, 3, 6, 0, 0, 109, 0),
    @(43996, 87598, 212, 1496, 1, 7, 1, 0, 109, 0),
    @(43997, 88165, 567, 1499, 3, 7, 1, 0, 109, 0),
    @(43998, 89151, 986, 1503, 4, 7, 1, 0, 109, 0),
    @(43999, 90103, 952, 1511, 8, 6, 1, 1, 109, 0),
    @(44000, 91005, 902, 1513, 2, 8, 1, 0, 109, 0),
    @(44001, 92152, 1147, 1519, 6, 6, 1, 2, 109, 0),
    @(44002, 92919, 758, 1520, 1, 6, 1, 2, 109, 0),
    @(44003, 93181, 271, 1521, 1, 6, 1, 0, 109, 0),
    @(44004, 94165, 984, 1534, 13, 5, 1, 1, 109, 0),
    @(44005, 95387, 1222, 1541, 7, 7, 2, 0, "111*", 0),
    @(44006, 96599, 1212, 1547, 6, 7, 2, 0, "111*", 0),
    @(44007, 97442, 843, 1558, 11, 8, 2, 0, "111*", 0),
    @(44008, 98320, 878, 1572, 14, 8, 1, 0, "111*", 0),
    @(44009, 98945, 625, 1581, 9, 7, 0, 1, "111*", 0),
    @(44010, 99245, 300, 1585, 4, 8, 0, 0, "111*", 0),
    @(44011, 100330, 1085, 1600, 15, 8, 0, 0, "111*", 0),
    @(44012, 101729, 1399, 1613, 13, 8, 0, 0, "111*", 0),
    @(44013, 102927, 1198, 1663, 21, 9, 0, 0, "111*", 0),
    @(44014, 104201, 1274, 1649, 16, 10, 0, 0, 111, 0)
)

$startRow = 76
$r = $startRow
foreach ($row in $data) {
    $prev = $r - 1
    $src = $ws.Range("A" + $prev + ":J" + $prev)
    $dst = $ws.Range("A" + $r + ":J" + $r)
    $src.Copy($dst)

    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }

    $r = $r + 1
}
$lastRow = $r - 1

# Grow the table / autofilter to cover the freshly written rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:J" + $lastRow))

# Match the author's final selection / scroll position.
$ws.Range("I" + $lastRow).Select()
